# Word COM-interop script: KM Airlines -> Swift (and Regional) Airlines
#
# The original document refers to "KM Airlines" in three spots:
#   1. The title line ("...for KM Airlines-Dulles VA")
#   2. The opening of the outline paragraph ("KM Airlines is a charter
#      airline company ...")
#   3. Mid-sentence in the same paragraph ("The CEO of KM Airlines has
#      requested ...")
#
# The edit renames the company to "Swift" (and tweaks the business
# description to "a charter and Regional airline company"). Word's real
# edit history shows the "_GoBack" bookmark (which always tracks the most
# recent edit point) moving from the title-line occurrence to the
# "CEO of Swift" occurrence, since that is where editing finished.
#
# Because collapsing/retyping text tends to make the interop host merge
# touching runs that share identical formatting, we rebuild each affected
# paragraph wholesale from OOXML (`Range.InsertXML`) so the exact run
# boundaries from the target document are preserved.

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits right after the first "KM" (title
# line). It is being relocated, so drop it now; it gets re-created at its
# new home when paragraph 7's XML is rebuilt below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Paragraph 1: title line ------------------------------------------------
$titleFind = $d.Content
$titleFind.Find.Text = "Project title:"
$null = $titleFind.Find.Execute()
$titlePara = $titleFind.Paragraphs(1).Range
$titleXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidR="008B5017" w:rsidRPr="00CB2623" w:rsidRDefault="00CB2623" w:rsidP="00CB2623"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00CB2623"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>P</w:t></w:r><w:r w:rsidR="008B5017" w:rsidRPr="00CB2623"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>roject title:</w:t></w:r><w:r w:rsidR="008B5017" w:rsidRPr="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Analyzing </w:t></w:r><w:r w:rsidR="008B5017" w:rsidRPr="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Turnover </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">for </w:t></w:r><w:r w:rsidR="009A08A7" w:rsidRPr="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Swift</w:t></w:r><w:r w:rsidR="008B5017" w:rsidRPr="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Airlines-Dulles VA</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.InsertXML($titleXml)

# --- Paragraph 7: outline / description paragraph --------------------------
$bodyFind = $d.Content
$bodyFind.Find.Text = "KM Airlines is a charter"
$null = $bodyFind.Find.Execute()
$bodyPara = $bodyFind.Paragraphs(1).Range
$bodyXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p w:rsidR="009A08A7" w:rsidRDefault="009A08A7" w:rsidP="009A08A7"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Swift</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Airlines is a charter </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">and Regional </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">airline company currently located in northern Virginia </w:t></w:r><w:r w:rsidR="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">area. It currently flies from DCA, IAD and BWI Airports. The IAD airport </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">of location </w:t></w:r><w:r w:rsidR="00C21662"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>is experiencing</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> a high turnover. The CEO of </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Swift</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Airlines has requested the financial analyst and HR analyst </w:t></w:r><w:r w:rsidR="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">of his organization </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">to provide </w:t></w:r><w:r w:rsidR="00CB2623"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>an analysis on the causes for the high turnover as well as recommendations.</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bodyPara.InsertXML($bodyXml)

Write-Output "Done. Title: $($d.Paragraphs(1).Range.Text)"
